# The commit adds one new daily price observation for
# "Feria Lagunitas de Puerto Montt - Ciboulette" dated 2022-10-11 (serial 44845,
# volume 240) into the weekly-ordered table. It is inserted as row 189,
# pushing the former rows 189:269 down to 190:270 (dimension grows from
# A1:R269 to A1:R270); no other existing data changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 189; Excel shifts rows 189:269 -> 190:270
# (and their formatting, e.g. the date style on column D) automatically.
$ws.Rows.Item(189).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(189, 1).Value = 4
$ws.Cells.Item(189, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(189, 3).Value = "Los Lagos"
$ws.Cells.Item(189, 4).Value = 44845
$ws.Cells.Item(189, 5).Value = 10
$ws.Cells.Item(189, 6).Value = 100112039
$ws.Cells.Item(189, 7).Value = "Ciboulette"
$ws.Cells.Item(189, 8).Value = "Sin especificar"
$ws.Cells.Item(189, 9).Value = "Primera"
$ws.Cells.Item(189, 10).Value = 240
$ws.Cells.Item(189, 11).Value = 2500
$ws.Cells.Item(189, 12).Value = 2500
$ws.Cells.Item(189, 13).Value = 2500
$ws.Cells.Item(189, 14).Value = "$/docena de atados"
$ws.Cells.Item(189, 15).Value = "Región Metropolitana"
$ws.Cells.Item(189, 16).Value = 833
$ws.Cells.Item(189, 17).Value = 3
$ws.Cells.Item(189, 18).Value = "Hortaliza"
